$wb = $excel.ActiveWorkbook

# --- 1. Update status text "Ready for handoff" -> "In Translation" ---
# The "Overview" sheet carries the per-language status in columns E (zh-cn) and
# F (de-de) for each content row; the per-language sheets ("zh-cn"/"de-de") carry
# the same status in column C ("Status"). Update every cell that shows the old
# status so the shared string collapses/rewrites cleanly, matching all rows.
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewUsed = $wsOverview.UsedRange
for ($r = 1; $r -le $overviewUsed.Rows.Count; $r++) {
    foreach ($c in 5, 6) {
        $cell = $wsOverview.Cells.Item($r, $c)
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

foreach ($sheetName in "zh-cn", "de-de") {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Shrink the status columns now that "In Translation" is shorter than
#        "Ready for handoff" (re-generated/auto-fit report widths) ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wb.Worksheets.Item("zh-cn").Columns.Item(3).ColumnWidth = 12.5
$wb.Worksheets.Item("de-de").Columns.Item(3).ColumnWidth = 12.5
